$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force A19 to be stored as literal text (not auto-converted to a date
# serial number) by temporarily marking it as a Text-formatted cell before
# assigning the value, matching how the source data was produced.
$ws.Range("A19").NumberFormat = "@"
$ws.Range("A19").Value = "2026-02-05"

# "21:38:36" is not parsed as a date/time value, it stays text naturally.
$ws.Range("B19").Value = "21:38:36"

# Restore the original row style (s="2", Times New Roman / General format)
# on A19 and B19 by copying formatting from the row above, since new cells
# default to the workbook's base style otherwise.
$ws.Range("A18:B18").Copy()
$ws.Range("A19").PasteSpecial(-4122)

# Clear the old "Total Duration:" label and its computed value; the cells
# keep their existing style (s="2") but hold no content any more.
$ws.Range("C19").Value = $null
$ws.Range("D19").Value = $null

$excel.CutCopyMode = 0
